# Generate Report for handoff
# Update the "Latest Handoff Datetime" (column D, row 5) for the
# 6715fa04-... file in both the zh-cn and de-de handoff status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-28 10:38:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-28 10:38:49"
